$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Green = V (correct), Red = X (incorrect) — same convention used throughout the sheet
$GREEN = 32768   # RGB(0,128,0) -> 00008000
$RED   = 255     # RGB(255,0,0) -> 00FF0000

# row | predictedScore (D) | correct (F: V or X)
$data = @(
    ,@(2, "1-1", "V")
    ,@(3, "0-2", "X")
    ,@(4, "1-1", "X")
    ,@(5, "1-1", "V")
    ,@(6, "2-1", "V")
    ,@(7, "2-0", "X")
    ,@(8, "1-1", "X")
    ,@(9, "1-2", "V")
    ,@(10, "1-1", "X")
    ,@(11, "0-2", "X")
    ,@(12, "3-1", "V")
    ,@(13, "1-1", "X")
    ,@(14, "1-1", "X")
    ,@(15, "1-2", "X")
    ,@(16, "1-2", "V")
    ,@(17, "1-1", "X")
    ,@(18, "1-1", "X")
    ,@(19, "2-1", "V")
    ,@(20, "2-1", "V")
    ,@(21, "2-1", "X")
    ,@(22, "1-2", "V")
    ,@(23, "1-1", "V")
    ,@(24, "1-1", "X")
    ,@(25, "1-1", "X")
    ,@(26, "1-1", "V")
    ,@(27, "2-1", "V")
    ,@(28, "1-1", "V")
    ,@(29, "1-1", "V")
    ,@(30, "3-1", "X")
    ,@(31, "1-2", "V")
    ,@(32, "1-1", "V")
    ,@(33, "1-2", "V")
    ,@(34, "1-1", "V")
    ,@(35, "1-1", "X")
    ,@(36, "2-0", "V")
    ,@(37, "2-1", "V")
    ,@(38, "1-1", "X")
    ,@(39, "1-2", "X")
    ,@(40, "1-3", "V")
    ,@(41, "1-2", "V")
    ,@(42, "1-1", "V")
    ,@(43, "1-2", "V")
    ,@(44, "1-1", "X")
    ,@(45, "0-1", "V")
    ,@(46, "2-1", "V")
    ,@(47, "0-2", "V")
    ,@(48, "1-2", "V")
    ,@(49, "2-1", "X")
    ,@(50, "1-2", "V")
    ,@(51, "2-1", "X")
    ,@(52, "1-1", "X")
    ,@(53, "1-1", "V")
    ,@(54, "1-1", "X")
    ,@(55, "2-1", "X")
    ,@(56, "2-1", "X")
    ,@(57, "1-1", "V")
    ,@(58, "2-1", "V")
    ,@(59, "0-2", "V")
    ,@(60, "1-1", "X")
    ,@(61, "1-2", "V")
    ,@(62, "1-2", "X")
    ,@(63, "1-2", "V")
    ,@(64, "1-1", "X")
    ,@(65, "1-1", "V")
    ,@(66, "1-2", "V")
    ,@(67, "2-2", "X")
    ,@(68, "0-2", "V")
)

foreach ($entry in $data) {
    $row = $entry[0]
    $predicted = $entry[1]
    $correct = $entry[2]

    $ws.Cells.Item($row, 4).Value = $predicted
    $ws.Cells.Item($row, 6).Value = $correct

    if ($correct -eq "V") {
        $ws.Cells.Item($row, 6).Interior.Color = $GREEN
    } else {
        $ws.Cells.Item($row, 6).Interior.Color = $RED
    }
}
